$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/eff597e02abd8c76d65af301c806c02429ee2ef3/e2e/"
$fileName = "e6803c1e-b5e5-408b-a400-cfb28f92e81e.md"
$fileUrl = $baseUrl + $fileName

# ---------------------------------------------------------------
# Sheet 1: Overview  (new row 3)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A3").Value = $fileName
$ws1.Range("B3").Value = "e2e\" + $fileName
$ws1.Range("C3").Value = ".md"
$ws1.Range("D3").Value = ""
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-12 22:50:25"
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Hyperlinks.Add($ws1.Range("B3"), $fileUrl, "", "", "e2e\" + $fileName) | Out-Null

$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:G3"))

# ---------------------------------------------------------------
# Sheet 2: zh-cn  (new row 3)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A3").Value = $fileName
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "'False"
$ws2.Range("G3").Value = "e6803c1e-b5e5-408b-a400-cfb28f92e81e.5b62383db7a32374694724733a3275c2bf01b44c.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-12 22:50:18"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = ""
$ws2.Range("K3").Value = "0001-01-01 00:00:00"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("L3").Value = ""
$ws2.Range("M3").Value = "'True"
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "'False"
$ws2.Range("P3").Value = ""

$ws2.Hyperlinks.Add($ws2.Range("A3"), $fileUrl, "", "", $fileName) | Out-Null

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------
# Sheet 3: de-de  (new row 3)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A3").Value = $fileName
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "'False"
$ws3.Range("G3").Value = "e6803c1e-b5e5-408b-a400-cfb28f92e81e.5b62383db7a32374694724733a3275c2bf01b44c.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-12 22:50:25"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = ""
$ws3.Range("K3").Value = "0001-01-01 00:00:00"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("L3").Value = ""
$ws3.Range("M3").Value = "'True"
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "'False"
$ws3.Range("P3").Value = ""

$ws3.Hyperlinks.Add($ws3.Range("A3"), $fileUrl, "", "", $fileName) | Out-Null

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:P3"))

Write-Output "done"
